$wb = $excel.ActiveWorkbook

# --- Sheet "leaderboard2" ---
$ws2 = $wb.Worksheets.Item("leaderboard2")
$ws2.Range("C3").Value = "BKZRackham"
$ws2.Range("D3").Value = 23
$ws2.Range("C4").Value = "ArtyumsM"
$ws2.Range("D4").Value = 9
$ws2.Range("B13").Value = "Dernière update le 26.03.25 à 00:34"

# --- Sheet "leaderboard3" ---
$ws3 = $wb.Worksheets.Item("leaderboard3")
$ws3.Range("C3").Value = "ArtyumsM"
$ws3.Range("C4").Value = "Lokys"
$ws3.Range("C5").Value = "Machoppeur_"
$ws3.Range("C6").Value = "BKZRackham"
$ws3.Range("B13").Value = "Dernière update le 26.03.25 à 00:34"
